$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.994.87"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "2.737.11"
$ws.Range("E3").Value = "  +3.58%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.08"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.14"
$ws.Range("E6").Value = "  +6.60%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.548"
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").Value = "2.735.56"
$ws.Range("E9").Value = "  +3.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.147"
$ws.Range("E10").Value = "  +3.45%  "
$ws.Range("E11").Value = "  +5.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.35"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.79"
$ws.Range("E14").Value = "  +3.27%  "
$ws.Range("D15").Value = "3.234.14"
$ws.Range("E15").Value = "  +3.50%  "
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("D17").Value = "68.837.29"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "2.727.60"
$ws.Range("E18").Value = "  +3.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.98"
$ws.Range("E19").Value = "  +5.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "374.52"
$ws.Range("E20").Value = "  +4.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.69"
$ws.Range("E21").Value = "  +4.99%  "
$ws.Range("E22").Value = "  +3.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.99"
$ws.Range("E23").Value = "  +5.51%  "
$ws.Range("E24").Value = "  +3.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.86"
$ws.Range("E25").Value = "  -1.58%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  +3.72%  "
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("E29").Value = "  +3.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "586.65"
$ws.Range("E30").Value = "  +5.30%  "
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.35"
$ws.Range("E32").Value = "  +4.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.46"
$ws.Range("E33").Value = "  +5.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.97"
$ws.Range("E34").Value = "  +5.67%  "
$ws.Range("E35").Value = "  +4.66%  "
$ws.Range("E36").Value = "  +4.90%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.06"
$ws.Range("E38").Value = "  +1.96%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "161.22"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.383"
$ws.Range("E40").Value = "  +3.44%  "
$ws.Range("E41").Value = "  +4.40%  "
$ws.Range("E42").Value = "  +3.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "18.00"
$ws.Range("E43").Value = "  +1.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.67"
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("D46").Value = "0.0₆0314"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.06"
$ws.Range("E47").Value = "  +2.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.55"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E50").Value = "  +7.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.607"
$ws.Range("E51").Value = "  +7.43%  "
